$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 53, pushing the existing
# rows 53-56 down to become rows 56-59.
$ws.Rows("53:55").Insert()

# New row 53: Especial, volumen 48
$ws.Range("A53").Value = 3
$ws.Range("B53").Value = "Femacal de La Calera"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44474
$ws.Range("E53").Value = 5
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107002
$ws.Range("J53").Value = "Chirimoya"
$ws.Range("K53").Value = "Cultivar IV Región"
$ws.Range("L53").Value = "Especial"
$ws.Range("M53").Value = 48
$ws.Range("N53").Value = 30000
$ws.Range("O53").Value = 30000
$ws.Range("P53").Value = 30000
$ws.Range("Q53").Value = "$/bandeja 10 kilos"
$ws.Range("R53").Value = "Provincia del Elquí"
$ws.Range("S53").Value = 3000
$ws.Range("T53").Value = 10

# New row 54: Primera, volumen 47
$ws.Range("A54").Value = 3
$ws.Range("B54").Value = "Femacal de La Calera"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44474
$ws.Range("E54").Value = 5
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100107
$ws.Range("H54").Value = "Otros"
$ws.Range("I54").Value = 100107002
$ws.Range("J54").Value = "Chirimoya"
$ws.Range("K54").Value = "Cultivar IV Región"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 47
$ws.Range("N54").Value = 27000
$ws.Range("O54").Value = 27000
$ws.Range("P54").Value = 27000
$ws.Range("Q54").Value = "$/bandeja 10 kilos"
$ws.Range("R54").Value = "Provincia del Elquí"
$ws.Range("S54").Value = 2700
$ws.Range("T54").Value = 10

# New row 55: Segunda, volumen 38
$ws.Range("A55").Value = 3
$ws.Range("B55").Value = "Femacal de La Calera"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44474
$ws.Range("E55").Value = 5
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100107
$ws.Range("H55").Value = "Otros"
$ws.Range("I55").Value = 100107002
$ws.Range("J55").Value = "Chirimoya"
$ws.Range("K55").Value = "Cultivar IV Región"
$ws.Range("L55").Value = "Segunda"
$ws.Range("M55").Value = 38
$ws.Range("N55").Value = 23000
$ws.Range("O55").Value = 23000
$ws.Range("P55").Value = 23000
$ws.Range("Q55").Value = "$/bandeja 10 kilos"
$ws.Range("R55").Value = "Provincia del Elquí"
$ws.Range("S55").Value = 2300
$ws.Range("T55").Value = 10

# Make sure the date cells keep the same date number format used by the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D53:D55").NumberFormat = $ws.Range("D56").NumberFormat
